# The deck had a blank, unused slide (Title/Content placeholders only,
# no text) sitting right before the final "Questions?" slide. Remove it.
#
# It carries p:sldId id="271" in the presentation's slide list (the last
# slide, "Questions?", keeps sldId 272 and simply shifts up to take its
# place). Locate it by SlideID rather than a hard-coded index so the
# script stays correct even if slide order assumptions drift.

$p = $ppt.ActivePresentation

$targetSlideId = 271
$slideToDelete = $null

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq $targetSlideId) {
        $slideToDelete = $candidate
        break
    }
}

if ($slideToDelete -ne $null) {
    $slideToDelete.Delete()
} else {
    # Fallback: the blank slide was the second-to-last one in the deck.
    $p.Slides.Item($p.Slides.Count - 1).Delete()
}
